$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.836.38"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.799.43"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.48"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.28"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.796.53"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  +12.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.80"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.435.33"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.781.48"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.890.14"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.10"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.67"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  -6.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.15"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.946.50"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.79"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.764.32"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.107"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.80"
$ws.Range("E38").Value = "  +12.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.93"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.84"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000291"
$ws.Range("E47").Value = "  +6.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "407.50"
$ws.Range("E48").Value = "  -5.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.48"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.55"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0357"
$ws.Range("E51").Value = "  +0.34%  "
